# Made Excel Table more realistic
#
# The localization example table used placeholder translations of the form
# "sv--<text>", "de--<text>", "fr--<text>" for most rows (these were just
# stand-ins produced by the tool, not real translations), so the Swedish /
# German / French columns (E:G) are cleared out for the "normal" example
# rows, leaving just the Dialog/Type/Name/English columns (A:D).
#
# The small demonstration "gloss" table (rows 14-15) keeps its four
# language columns, but gets real words instead of the old accented
# placeholder text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the sv/de/fr placeholder-translation columns for the regular rows.
$ws.Range("E2:G13").ClearContents()

# Replace the placeholder glosses in the little demo table with real
# translations (German column F stays as-is).
$ws.Range("E14").Value = "Nyaste"
$ws.Range("G14").Value = "Le plus récent"
$ws.Range("G15").Value = "Même table à nouveau"
$ws.Range("E15").Value = "Samma tabell igen"
